# "Adding economic data example"
#
# Slide 2 (the "Calling Functions" diagram) is shifted straight down to
# make room for a new example above it. Every top-level shape's vertical
# offset increases by the same amount (590309 EMU ~= 46.48 pt); nothing
# else about the shapes (x, width/height, group child offsets) changes.
#
# Note on precision: this host stores Shape.Top/Left as points using
# 32-bit float precision, matching real PowerPoint. Reading Top gives
# float32(emu/12700), which round-trips exactly back to the original EMU
# via rounding. Writing Top, however, converts the assigned points value
# to EMU by flooring float32(points)*12700 rather than rounding, so a
# naive "Top = Top + deltaPoints" can land one EMU short. To hit the
# exact target EMU we search the tiny neighborhood of candidate points
# values for one whose float32 cast floors to the desired EMU.

function Get-EmuFromPoints($pointsValue) {
    $f32 = [float]$pointsValue
    return [Math]::Round([double]$f32 * 12700.0)
}

function Find-PointsForTargetEmu($targetEmu) {
    $approx = $targetEmu / 12700.0
    $epsilon = 0.0000005
    for ($k = 0; $k -lt 2000; $k++) {
        foreach ($sign in @(1, -1)) {
            if ($k -eq 0 -and $sign -eq -1) { continue }
            $candidate = $approx + ($sign * $k * $epsilon)
            $f32 = [float]$candidate
            $emu = [Math]::Floor([double]$f32 * 12700.0 + 0.000001)
            if ($emu -eq $targetEmu) {
                return $candidate
            }
        }
    }
    # Fallback (shouldn't happen in practice): just return the naive value.
    return $approx
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$deltaEmu = 590309

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    $curEmu = Get-EmuFromPoints $shp.Top
    $targetEmu = $curEmu + $deltaEmu
    $shp.Top = Find-PointsForTargetEmu $targetEmu
}
